$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.40"
$ws.Range("E2").Value = "'0.78%"
$ws.Range("E3").Value = "'4.70%"
$ws.Range("D4").Value = "'5.651"
$ws.Range("E4").Value = "'-0.93%"
$ws.Range("E5").Value = "'3.05%"
$ws.Range("D6").Value = "'2.040"
$ws.Range("E6").Value = "'1.69%"
$ws.Range("D7").Value = "'8.767"
$ws.Range("E9").Value = "'0.54%"
$ws.Range("D10").Value = "'0.9253"
$ws.Range("E10").Value = "'0.21%"
$ws.Range("D11").Value = "'0.1274"
$ws.Range("E11").Value = "'1.28%"
$ws.Range("D12").Value = "'0.1962"
$ws.Range("E12").Value = "'0.02%"
$ws.Range("D13").Value = "'0.09389"
$ws.Range("E13").Value = "'1.75%"
$ws.Range("D14").Value = "'0.03926"
$ws.Range("E14").Value = "'10.11%"
$ws.Range("D15").Value = "'0.1060"
$ws.Range("E15").Value = "'0.88%"
$ws.Range("D16").Value = "'0.001305"
$ws.Range("E16").Value = "'0.19%"
$ws.Range("D17").Value = "'0.006153"
$ws.Range("E17").Value = "'-3.46%"
$ws.Range("D19").Value = "'3.438"
$ws.Range("E19").Value = "'2.20%"
$ws.Range("E20").Value = "'1.47%"
$ws.Range("D21").Value = "'8.359"
$ws.Range("E21").Value = "'-4.58%"
$ws.Range("E22").Value = "'1.85%"
$ws.Range("D24").Value = "'0.04403"
$ws.Range("E24").Value = "'0.12%"
$ws.Range("D25").Value = "'0.001257"
$ws.Range("E25").Value = "'-0.16%"
$ws.Range("D26").Value = "'0.004319"
$ws.Range("E26").Value = "'-6.49%"
$ws.Range("E27").Value = "'0.99%"
$ws.Range("D39").Value = "'0.02776"
$ws.Range("E39").Value = "'11.24%"
$ws.Range("D40").Value = "'0.05521"
$ws.Range("E40").Value = "'3.32%"
$ws.Range("D41").Value = "'0.007936"
$ws.Range("E41").Value = "'6.15%"
$ws.Range("D42").Value = "'0.1422"
$ws.Range("E42").Value = "'0.90%"
$ws.Range("D43").Value = "'0.008935"
$ws.Range("E43").Value = "'-9.80%"
$ws.Range("D44").Value = "'0.002142"
$ws.Range("E44").Value = "'1.33%"
$ws.Range("D45").Value = "'0.01189"
$ws.Range("E45").Value = "'20.07%"
$ws.Range("D46").Value = "'0.00007006"
$ws.Range("E46").Value = "'4.83%"
$ws.Range("E47").Value = "'0.15%"
$ws.Range("D48").Value = "'0.003191"
$ws.Range("E48").Value = "'5.01%"
$ws.Range("E49").Value = "'0.11%"
$ws.Range("E50").Value = "'0.15%"
$ws.Range("E51").Value = "'0.15%"
